$d = $word.ActiveDocument

# Locate the "Requisitos" detail paragraph ("LOB1004: Cálculo II (Requisito
# fraco)"). Immediately after it the document has three paragraphs that
# must be removed:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# The paragraph that follows those three (an empty paragraph right before
# the page-break paragraph) must be left untouched.
$anchorIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1004*Requisito fraco*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -ne $null) {
    $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

    $rng = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $rng.Delete()
}
